$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# batch_size (col J) 32 -> 8, negatives (col K) 16 -> 4 for the 5 run rows (2-6)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 10).Value = 8
    $ws.Cells.Item($r, 11).Value = 4
}

# Column K ("negatives") gets an explicit best-fit custom width now that its
# values are shorter, matching the sheet's post-edit autofit.
$ws.Columns.Item(11).ColumnWidth = 8.6

# Selection moved from C10 to K8
$ws.Range("K8").Select()
